# Applies the "Space" -> "Political Challenges" essay rewrite described by the
# commit diff. Uses Find/Replace for straightforward sentence-for-sentence
# swaps (these keep each sentence in its own run because every sentence in
# this document already lives in its own <w:r>), and Range-based surgery for
# the spots where whole runs are deleted, merged, or newly inserted.

$d = $word.ActiveDocument

# Find.Execute's built-in Replace (wdReplaceOne/All) runs the match through
# Word's AutoCorrect/AutoFormat pipeline, which silently turns straight
# apostrophes into curly ones. Locating the match with Find (no replacement,
# wdFindStop) and then assigning Range.Text ourselves sidesteps that, while
# still reproducing Word's normal "merge adjacent identically-formatted runs
# when their text changes" behaviour.
function Replace-Text($old, $new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output ("NOT FOUND: " + $old)
        return
    }
    $r.Text = $new
}

# ---------------------------------------------------------------------------
# Title / byline / contact block
# ---------------------------------------------------------------------------

Replace-Text "Unveiling the Mysteries of Space: A Journey Through the Cosmos" "Exploring the Enigma of Political Challenges in a Dynamic World"

Replace-Text "Evelyn Mitchell" "Eleanor Townsend"

# The email address is spread across 5 runs: "evelyn" / "." /
# "mitchell@astronomyenthusiast" / "." / "org". The target keeps 3 runs:
# "townsendeleano@anonymousemai1" / "." / "com" -- i.e. the first run's text
# changes, the "." after it is untouched, and the last three runs collapse
# into a single "com" run. Do the two edits as one pass each, outermost
# (longer) match first so the untouched middle "." is never part of either
# search string.
Replace-Text "evelyn" "townsendeleano@anonymousemai1"
Replace-Text "mitchell@astronomyenthusiast.org" "com"

# ---------------------------------------------------------------------------
# Opening paragraph (Introduction)
# ---------------------------------------------------------------------------

Replace-Text "Among the vast expanse of the cosmos, we, as inhabitants of Earth, are but minuscule entities captivated by the allure of the stars that twinkle above us" "In today's ever-evolving political landscape, understanding the challenges faced by nations is paramount"

Replace-Text " The universe, an infinite canvas adorned with celestial wonders, beckons us to embark on a journey of discovery, to unravel its secrets and comprehend our place within its boundless realms" " This essay delves into the intricacies of political hurdles, examining their causes, exploring their impact, and proposing potential solutions"

Replace-Text " From the celestial ballet of planets to the mesmerizing dance of distant galaxies, the cosmos whispers tales of cosmic evolution and the profound interconnectedness of all existence" " By examining historical and contemporary examples, this essay sheds light on the complex interplay between political systems, societal factors, economic conditions, and global dynamics"

Replace-Text "As we pierce through the veil of the unknown, unraveling the cosmic tapestry thread by thread, we encounter celestial bodies of captivating beauty and profound significance" "The political framework of a nation is a delicate tapestry woven from the threads of power, influence, and governance"

Replace-Text " The blazing suns, like miniature universes, ignite galaxies with their radiant energy, fueling the birth and evolution of celestial systems" " Various institutions, including governments, political parties, and electoral systems, shape the political landscape and determine how decisions are made"

# This replacement spans what used to be three runs ("Planets, trapped...",
# ".", " Moons, faithful companions...") and collapses them into one, exactly
# like the diff shows (the trailing "." run right after is left untouched).
Replace-Text " Planets, trapped in the gravitational embrace of their parent stars, revolve in intricate patterns, each harboring unique geological and atmospheric conditions. Moons, faithful companions, grace these planetary realms, adding celestial intrigue and enhancing the dynamic interplay of gravitational forces" " Understanding the dynamics of these institutions and their interactions is crucial in comprehending the challenges faced by nations"

Replace-Text "The realm of stars, however, extends beyond the familiar confines of our solar system" "Political challenges often arise from the tension between competing interests and ideologies within a society"

Replace-Text " Scattered across the vastness of space, stars of diverse sizes, colors, and life cycles ignite the celestial tapestry" " Differing perspectives on issues such as economic policies, social welfare, and foreign relations can lead to political gridlock and hinder progress"

Replace-Text " From massive blue giants, radiating intense heat and light, to aging red giants, gently shedding their outer layers, stars traverse a mesmerizing evolutionary path, painting the cosmos with a kaleidoscope of cosmic hues" " Additionally, societal factors like poverty, inequality, and discrimination can exacerbate political challenges, further complicating the path to resolution"

# ---------------------------------------------------------------------------
# The closing sentence of that paragraph ("As these stellar giants ...
# energy") is wiped out and replaced with a big new block of extra sentences
# (all still inside the same paragraph, separated by manual line breaks),
# finishing with "... can create a more stable and just world" -- the
# trailing "." run right after it is untouched, exactly as in the diff.
# ---------------------------------------------------------------------------

$d.Content.Find.Execute(" As these stellar giants reach the end of their luminous journeys, they may explode in spectacular supernovae, leaving behind remnants that challenge our understanding of matter and energy", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# Find the run that now holds nothing but the break, right before the
# trailing "." that ends the paragraph, and insert the new material after it.
$anchorRange = $d.Content
$anchorRange.Find.Execute("Summary", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$summaryStart = $anchorRange.Start

$insertPoint = $d.Range($summaryStart - 2, $summaryStart - 2)

$newBody = "Body:" + "`v" + "`v" + "Political systems face myriad challenges, ranging from internal conflicts and power struggles to external pressures and global uncertainties" + "." + " Internal conflicts, often rooted in historical grievances or ethnic tensions, can lead to political instability and violence" + "." + " Power struggles between different factions within a government can paralyze decision-making and hinder effective governance" + "." + " External pressures, such as economic sanctions or military threats, can also pose significant challenges to a nation's political stability" + "." + "`v" + "`v" + "Navigating the complexities of global dynamics presents another set of challenges for nations" + "." + " In an increasingly interconnected world, actions taken by one nation can have far-reaching consequences for others" + "." + " Globalization, while fostering economic growth and cultural exchange, also amplifies the interconnectedness of political systems" + "." + " This interdependence can lead to diplomatic disputes, trade conflicts, and geopolitical tensions, requiring skillful diplomacy and cooperation to resolve" + "." + "`v" + "`v" + "To address these challenges, nations must seek collaborative solutions that prioritize dialogue, negotiation, and compromise" + "." + " Strengthening democratic institutions, promoting transparency and accountability, and investing in education and public awareness can help build a more informed and engaged citizenry capable of holding their leaders accountable" + "." + " Additionally, fostering international cooperation, promoting economic development, and addressing global issues collectively can create a more stable and just world"

$insertPoint.InsertBefore("`v" + $newBody)

# ---------------------------------------------------------------------------
# Summary paragraph
# ---------------------------------------------------------------------------

Replace-Text "Our journey through the cosmos unveils a universe teeming with wonder and mystery" "Political challenges are multifaceted and ever-changing, influenced by a complex interplay of internal dynamics, societal factors, and global forces"

Replace-Text " From the celestial ballet of planets to the majestic tapestry of stars, the universe beckons us to explore its enigmatic depths and contemplate the profound interconnectedness of all existence" " Understanding the intricacies of these challenges is essential for developing effective solutions"

Replace-Text " As we continue to unravel the secrets of space, we deepen our understanding of our place within the boundless realms of the universe, inspiring awe, curiosity, and a profound appreciation for the intricate beauty of creation" " By promoting dialogue, collaboration, and international cooperation, nations can navigate the complexities of political challenges and strive for a more stable, just, and peaceful world"

# ---------------------------------------------------------------------------
# Trailing empty paragraph added after the Summary paragraph, before sectPr.
# ---------------------------------------------------------------------------

$endPoint = $d.Range($d.Content.End, $d.Content.End)
$endPoint.InsertParagraphAfter() | Out-Null
